$d = $word.ActiveDocument

$d.Content.Find.Execute("60×14=840", $true, $false, $false, $false, $false, $true, 1, $false, "39×20=780", 2) | Out-Null
$d.Content.Find.Execute("26×13=338", $true, $false, $false, $false, $false, $true, 1, $false, "48×90=4320", 2) | Out-Null
$d.Content.Find.Execute("35×25=875", $true, $false, $false, $false, $false, $true, 1, $false, "80×25=2000", 2) | Out-Null
$d.Content.Find.Execute("32×41=1312", $true, $false, $false, $false, $false, $true, 1, $false, "54×26=1404", 2) | Out-Null
$d.Content.Find.Execute("35×50=1750", $true, $false, $false, $false, $false, $true, 1, $false, "44×64=2816", 2) | Out-Null
$d.Content.Find.Execute("18×73=1314", $true, $false, $false, $false, $false, $true, 1, $false, "54×90=4860", 2) | Out-Null
$d.Content.Find.Execute("79×64=5056", $true, $false, $false, $false, $false, $true, 1, $false, "60×67=4020", 2) | Out-Null
$d.Content.Find.Execute("28×71=1988", $true, $false, $false, $false, $false, $true, 1, $false, "32×43=1376", 2) | Out-Null
$d.Content.Find.Execute("14×55=770", $true, $false, $false, $false, $false, $true, 1, $false, "94×18=1692", 2) | Out-Null
$d.Content.Find.Execute("18×60=1080", $true, $false, $false, $false, $false, $true, 1, $false, "73×22=1606", 2) | Out-Null
$d.Content.Find.Execute("56×88=4928", $true, $false, $false, $false, $false, $true, 1, $false, "27×15=405", 2) | Out-Null
$d.Content.Find.Execute("61×95=5795", $true, $false, $false, $false, $false, $true, 1, $false, "14×99=1386", 2) | Out-Null
$d.Content.Find.Execute("89×92=8188", $true, $false, $false, $false, $false, $true, 1, $false, "62×87=5394", 2) | Out-Null
$d.Content.Find.Execute("74×34=2516", $true, $false, $false, $false, $false, $true, 1, $false, "22×81=1782", 2) | Out-Null
$d.Content.Find.Execute("73×59=4307", $true, $false, $false, $false, $false, $true, 1, $false, "22×51=1122", 2) | Out-Null
$d.Content.Find.Execute("88×80=7040", $true, $false, $false, $false, $false, $true, 1, $false, "48×23=1104", 2) | Out-Null
$d.Content.Find.Execute("20×32=640", $true, $false, $false, $false, $false, $true, 1, $false, "53×65=3445", 2) | Out-Null
$d.Content.Find.Execute("23×96=2208", $true, $false, $false, $false, $false, $true, 1, $false, "46×23=1058", 2) | Out-Null
$d.Content.Find.Execute("37×71=2627", $true, $false, $false, $false, $false, $true, 1, $false, "81×60=4860", 2) | Out-Null
$d.Content.Find.Execute("83×36=2988", $true, $false, $false, $false, $false, $true, 1, $false, "21×88=1848", 2) | Out-Null
$d.Content.Find.Execute("73×51=3723", $true, $false, $false, $false, $false, $true, 1, $false, "13×57=741", 2) | Out-Null
$d.Content.Find.Execute("44×37=1628", $true, $false, $false, $false, $false, $true, 1, $false, "29×89=2581", 2) | Out-Null
$d.Content.Find.Execute("87×82=7134", $true, $false, $false, $false, $false, $true, 1, $false, "39×51=1989", 2) | Out-Null
$d.Content.Find.Execute("12×76=912", $true, $false, $false, $false, $false, $true, 1, $false, "64×13=832", 2) | Out-Null
$d.Content.Find.Execute("85×85=7225", $true, $false, $false, $false, $false, $true, 1, $false, "63×73=4599", 2) | Out-Null
